# Apply stagnation-threshold (2n) / 3-neighbourhood results update to the
# "My Approach" column (E) on the "04 Sep" sheet. Column J recalculates
# automatically via its existing shared formula (=G-E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("04 Sep")

# row -> new "My Approach" (column E) value
$updates = @{
    8  = 80      # CON6
    9  = 124     # CON8
    10 = 184     # CON10
    11 = 263     # CON12
    12 = 348     # CON14
    15 = 42211   # NL6
    16 = 70237   # NL8
    17 = 133682  # NL10
    18 = 247736  # NL12
    19 = 350399  # NL14
    20 = 66      # CIRC6
    22 = 296     # CIRC10
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}

$wb.Application.CalculateFull()
